# continue cost reduction improvements
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- formula / value updates (cost reduction pass) ---
$ws.Range("B29").Formula = "=B87-B1-B9"
$ws.Range("B87").Formula = "=(B24-B64)/2+B8+B2+0.1"
$ws.Range("B89").Value = 20
$ws.Range("B91").Formula = "=FLOOR(B88-B89-B90-10,1)"

# --- view/selection updates ---
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("B30").Select()
$win.Zoom = 100
